$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.361.15'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = '2.095.84'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'251.95"
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'54.06"
$ws.Range("E8").Value = '  +20.41%  '
$ws.Range("D9").Value = "'62.48"
$ws.Range("E9").Value = '  +2.95%  '
$ws.Range("D10").Value = "'0.380"
$ws.Range("E10").Value = '  +4.26%  '
$ws.Range("D11").Value = "'0.0754"
$ws.Range("E11").Value = '  +4.96%  '
$ws.Range("E12").Value = '  +7.97%  '
$ws.Range("D13").Value = "'15.21"
$ws.Range("E13").Value = '  +5.65%  '
$ws.Range("D14").Value = '2.402.32'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = "'0.837"
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '2.100.03'
$ws.Range("D17").Value = "'5.24"
$ws.Range("E17").Value = '  +7.09%  '
$ws.Range("D18").Value = '37.312.73'
$ws.Range("E18").Value = '  +2.08%  '
$ws.Range("D19").Value = "'73.21"
$ws.Range("E19").Value = '  +2.82%  '
$ws.Range("D20").Value = "'14.52"
$ws.Range("E20").Value = '  +15.25%  '
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  +5.34%  '
$ws.Range("D22").Value = "'241.51"
$ws.Range("E23").Value = '  +7.59%  '
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").Value = "'171.38"
$ws.Range("E26").Value = '  +1.40%  '
$ws.Range("E27").Value = '  +5.27%  '
$ws.Range("D28").Value = "'20.93"
$ws.Range("E28").Value = '  +3.52%  '
$ws.Range("E29").Value = '  +5.38%  '
$ws.Range("E30").Value = '  +2.55%  '
$ws.Range("D31").Value = "'23.37"
$ws.Range("E31").Value = '  +8.62%  '
$ws.Range("E32").Value = '  +24.22%  '
$ws.Range("E33").Value = '  +4.90%  '
$ws.Range("D34").Value = "'0.0627"
$ws.Range("E34").Value = '  +8.09%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +6.60%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("E39").Value = '  -3.09%  '
$ws.Range("D40").Value = "'1.37"
$ws.Range("E40").Value = '  +3.77%  '
$ws.Range("D41").Value = "'4.95"
$ws.Range("E41").Value = '  +147.84%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = "'0.0229"
$ws.Range("E42").Value = '  +7.07%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'18.05"
$ws.Range("E43").Value = '  +13.80%  '
$ws.Range("E44").Value = '  +5.00%  '
$ws.Range("D45").Value = "'99.38"
$ws.Range("E45").Value = '  +2.90%  '
$ws.Range("D46").Value = "'0.0962"
$ws.Range("E46").Value = '  +17.61%  '
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = '1.335.19'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("E49").Value = '  +4.86%  '
$ws.Range("E50").Value = '  +8.07%  '
$ws.Range("E51").Value = '  +13.42%  '
